$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 24 with the next day's data
$ws.Range("A24").Value = 43928
$ws.Range("A24").NumberFormat = "YYYY\-MM\-DD"

$ws.Range("B24").Value = 586

# Update active selection to C24
$ws.Range("C24").Select()
